$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- D3: update E400 DS3 fit numbers (J0/J1) ---
$ws.Range("D3").Value = "`nE400 DS3`nJ0: 27.346 +/- 0.630`nJ1:  6.549 +/- 0.626;"

# --- E3: append new Ei400 DS3 fit block after the existing Ei400 numbers ---
$ws.Range("E3").Value = "Ei400`nJ0: 62.344 +/- 0.754; `nJ1: -5.894 +/- 0.854;`nJ2: -7.533 +/- 0.216`nEi400 DS3`nJ0: 65.663 +/- 3.055; `nJ1: -7.930 +/- 2.153; `nJ2: -8.047 +/- 0.476;"

# --- E2: rename Ei400 -> Ei400 DS2 and append the new E400 DS3 block ---
$ws.Range("E2").Value = "E200 NoA: `nJ0: 28.146 +/- 0.068;`nJ1: 12.931 +/- 0.064;`nJ2: -4.232 +/- 0.020 `nEi400 DS2`nJ0: 91.488 +/- 2.931; `nJ1: -11.325 +/- 2.388;`nJ2: -15.305 +/- 0.563`nE400 DS3`nJ0:64.120 +/- 3.159; `nJ1: -8.249 +/- 2.186;`nJ2: -7.536 +/- 0.524"

# --- C2: prepend a new "E200 DS_enhanced" / "J0: (running)" rich-text header ---
$ws.Range("C2").Value = "E200 DS_enhanced `nJ0: (running)`nEi400 DS3`n J0: 33.72 +/- 0.1`n"
$c2 = $ws.Range("C2").Characters(19, 13)
$c2.Font.Color = 5287936
$c2.Font.Name = "Calibri"
$c2.Font.Size = 11

# --- row 2 height grows to fit the extra line in C2 ---
$ws.Rows.Item(2).RowHeight = 210

# --- selection moves from D2 to F2 ---
$ws.Range("F2").Select()
